$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 12.517083
$ws.Range("D3").Value = 54.927585
